$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 34
$ws.Range("H34").Value = 3849.2727
$ws.Range("I34").Value = 3234.2
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 3234.2
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -3031.2
$ws.Range("N34").Value = -10406
# Row 36
$ws.Range("H36").Value = 3849.2727
$ws.Range("I36").Value = 3234.2
$ws.Range("J36").Value = 10000
$ws.Range("K36").Value = 3234.2
$ws.Range("L36").Value = 10000
$ws.Range("M36").Value = -2519.2
$ws.Range("N36").Value = -11430
# Row 55
$ws.Range("H55").Value = 947355
$ws.Range("I55").Value = 677.5
$ws.Range("J55").Value = 1838345.6
$ws.Range("K55").Value = 677.5
$ws.Range("L55").Value = 1838345.6
$ws.Range("M55").Value = -463.5
$ws.Range("N55").Value = -1838773.6
# Row 61
$ws.Range("H61").Value = 360
# Row 64
$ws.Range("H64").Value = 3560
$ws.Range("I64").Value = 3266.6667
$ws.Range("K64").Value = 3266.6667
$ws.Range("M64").Value = -3018.6667
# Row 67
$ws.Range("H67").Value = 3560
$ws.Range("I67").Value = 3266.6667
$ws.Range("K67").Value = 3266.6667
$ws.Range("M67").Value = -2408.6667
# Row 80
$ws.Range("H80").Value = 351
$ws.Range("I80").Value = 351
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1053
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -55
# Row 83
$ws.Range("H83").Value = 351
$ws.Range("I83").Value = 351
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 3159
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = 1833
# Row 86
$ws.Range("H86").Value = 961815
$ws.Range("I86").Value = 10625
$ws.Range("J86").Value = 2184773.5
$ws.Range("K86").Value = 10625
$ws.Range("L86").Value = 2184773.5
$ws.Range("M86").Value = -9502
$ws.Range("N86").Value = -2187019.5
# Row 89
$ws.Range("H89").Value = 961815
$ws.Range("I89").Value = 10625
$ws.Range("J89").Value = 2184773.5
$ws.Range("K89").Value = 53125
$ws.Range("L89").Value = 10923867.5
$ws.Range("M89").Value = -47509
$ws.Range("N89").Value = -10935099.5
# Row 132
$ws.Range("H132").Value = 6670508
$ws.Range("I132").Value = 8478017
$ws.Range("K132").Value = 25434051
$ws.Range("M132").Value = -25431521

$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 3717.0667
$ws.Range("I63").Value = 4836.5713
$ws.Range("J63").Value = 2737.5
$ws.Range("K63").Value = 4836.5713
$ws.Range("L63").Value = 2737.5
$ws.Range("M63").Value = -4150.5713
$ws.Range("N63").Value = -4109.5
# Row 66
$ws.Range("H66").Value = 3717.0667
$ws.Range("I66").Value = 4836.5713
$ws.Range("J66").Value = 2737.5
$ws.Range("K66").Value = 24182.8565
$ws.Range("L66").Value = 13687.5
$ws.Range("M66").Value = -20750.8565
$ws.Range("N66").Value = -20551.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1815.7174
$ws.Range("I20").Value = 1832.3125
$ws.Range("J20").Value = 1777.7858
$ws.Range("K20").Value = 1832.3125
$ws.Range("L20").Value = 1777.7858
$ws.Range("M20").Value = -1585.3125
$ws.Range("N20").Value = -2271.7858
# Row 23
$ws.Range("H23").Value = 7400
$ws.Range("J23").Value = 7400
$ws.Range("L23").Value = 7400
$ws.Range("N23").Value = -7966
# Row 94
$ws.Range("H94").Value = 705.8
$ws.Range("I94").Value = 699.7619
$ws.Range("J94").Value = 737.5
$ws.Range("K94").Value = 699.7619
$ws.Range("L94").Value = 737.5
$ws.Range("M94").Value = -248.7619
$ws.Range("N94").Value = -1639.5

$ws = $wb.Worksheets.Item("CRP")
# Row 55
$ws.Range("H55").Value = 2557.6667
$ws.Range("I55").Value = 936.5
$ws.Range("J55").Value = 5800
$ws.Range("K55").Value = 936.5
$ws.Range("L55").Value = 5800
$ws.Range("M55").Value = -621.5
$ws.Range("N55").Value = -6430

$ws = $wb.Worksheets.Item("CUL")
# Row 22
$ws.Range("H22").Value = 3710058.2
$ws.Range("I22").Value = 50000750
$ws.Range("J22").Value = 6802.96
$ws.Range("K22").Value = 150002250
$ws.Range("L22").Value = 20408.88
$ws.Range("M22").Value = -150002081
$ws.Range("N22").Value = -20746.88
# Row 27
$ws.Range("H27").Value = 3710058.2
$ws.Range("I27").Value = 50000750
$ws.Range("J27").Value = 6802.96
$ws.Range("K27").Value = 150002250
$ws.Range("L27").Value = 20408.88
$ws.Range("M27").Value = -150002148
$ws.Range("N27").Value = -20612.88
# Row 94
$ws.Range("H94").Value = 2238.8
$ws.Range("I94").Value = 1497
$ws.Range("J94").Value = 2733.3333
$ws.Range("K94").Value = 4491
$ws.Range("L94").Value = 8199.999899999999
$ws.Range("M94").Value = -3815
$ws.Range("N94").Value = -9551.999899999999
# Row 131
$ws.Range("H131").Value = 839.84
$ws.Range("J131").Value = 865.3117999999999
$ws.Range("L131").Value = 2595.9354
$ws.Range("N131").Value = -12675.9354

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 600.75
$ws.Range("I3").Value = 467.66666
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 467.66666
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -351.66666
$ws.Range("N3").Value = -1232

$ws = $wb.Worksheets.Item("LTW")
# Row 48
$ws.Range("H48").Value = 9952.666999999999
$ws.Range("I48").Value = 5000
$ws.Range("J48").Value = 10943.2
$ws.Range("K48").Value = 5000
$ws.Range("L48").Value = 10943.2
$ws.Range("M48").Value = -4339
$ws.Range("N48").Value = -12265.2
# Row 68
$ws.Range("H68").Value = 1504.4
$ws.Range("I68").Value = 1469
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1469
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -720
$ws.Range("N68").Value = -3498
# Row 71
$ws.Range("H71").Value = 1504.4
$ws.Range("I71").Value = 1469
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 7345
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -3601
$ws.Range("N71").Value = -17488
# Row 132
$ws.Range("H132").Value = 25007398
$ws.Range("I132").Value = 55557572
$ws.Range("J132").Value = 11800.363
$ws.Range("K132").Value = 166672716
$ws.Range("L132").Value = 35401.089
$ws.Range("M132").Value = -166670186
$ws.Range("N132").Value = -40461.089

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 850
$ws.Range("I126").Value = 375
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 1125
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = 1345
$ws.Range("N126").Value = -13190
# Row 132
$ws.Range("H132").Value = 11655582
$ws.Range("I132").Value = 35579.2
$ws.Range("J132").Value = 38470972
$ws.Range("K132").Value = 106737.6
$ws.Range("L132").Value = 115412916
$ws.Range("M132").Value = -104207.6
$ws.Range("N132").Value = -115417976
